$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E7").Value = 34
$ws.Range("G7").Value = 4
$ws.Range("H7").Value = 26

$ws.Range("E9").Value = 10
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 6

$ws.Range("E10").Value = 625

$ws.Range("E11").Value = 408
$ws.Range("F11").Value = 223
$ws.Range("H11").Value = 287

$ws.Range("E12").Value = 621
$ws.Range("F12").Value = 351
$ws.Range("H12").Value = 437

$ws.Range("E15").Value = 186

$ws.Range("E16").Value = 218
$ws.Range("F16").Value = 113
$ws.Range("H16").Value = 161

$ws.Range("E23").Value = 213

$ws.Range("E24").Value = 235

$ws.Range("E27").Value = 355
$ws.Range("F27").Value = 188
$ws.Range("H27").Value = 269

$ws.Range("E30").Value = 232

$ws.Range("E32").Value = 195
$ws.Range("F32").Value = 119
$ws.Range("H32").Value = 157

$ws.Range("E33").Value = 311

$ws.Range("E42").Value = 415

$ws.Range("E44").Value = 333
$ws.Range("F44").Value = 171
$ws.Range("H44").Value = 239

$ws.Range("E45").Value = 162

$ws.Range("E46").Value = 357

$ws.Range("E47").Value = 501
$ws.Range("F47").Value = 265
$ws.Range("H47").Value = 357

$ws.Range("E48").Value = 238

$ws.Range("E49").Value = 307

$ws.Range("E50").Value = 258

$ws.Range("F51").Value = 119
$ws.Range("H51").Value = 193

$ws.Range("E52").Value = 30
